$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: several row-pairs in the sheet had their match data (columns F..V,
# i.e. home team through match URL) swapped between two adjacent rows while
# the leading columns (A index, B country, C tournament, D season, E match
# date) stayed put. Swap the F:V payload back between each pair.
# ---------------------------------------------------------------------------

function Swap-MatchData($row1, $row2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $cell1 = $ws.Range($col + $row1)
        $cell2 = $ws.Range($col + $row2)
        $val1 = $cell1.Value()
        $val2 = $cell2.Value()
        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

$rowPairs = @(
    @(24, 25),
    @(36, 37),
    @(78, 79),
    @(86, 87),
    @(164, 165),
    @(167, 168),
    @(169, 170),
    @(172, 173),
    @(177, 178),
    @(179, 180)
)

foreach ($pair in $rowPairs) {
    Swap-MatchData $pair[0] $pair[1]
}

# ---------------------------------------------------------------------------
# Part 2: three new match rows were appended at the bottom of the sheet
# (rows 187-189, Indice 186-188). Copy the formatting (styles) of the last
# existing row down into the new rows, then populate their values.
# ---------------------------------------------------------------------------

$ws.Range("A186:V186").Copy()
$ws.Range("A187:V189").PasteSpecial(-4122)

$ws.Range("A187").Value = 186
$ws.Range("B187").Value = "italy"
$ws.Range("C187").Value = "serie-a"
$ws.Range("D187").Value = "2023-2024"
$ws.Range("E187").Value = 45298.52083333334
$ws.Range("F187").Value = "Empoli"
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = "AC Milan"
$ws.Range("I187").Value = 3
$ws.Range("J187").Value = 4.51
$ws.Range("K187").Value = "23/12/2024 23:02"
$ws.Range("L187").Value = 5.6
$ws.Range("M187").Value = "07/01/2024 12:21"
$ws.Range("N187").Value = 3.88
$ws.Range("O187").Value = "23/12/2024 23:02"
$ws.Range("P187").Value = 4.1
$ws.Range("Q187").Value = "07/01/2024 12:21"
$ws.Range("R187").Value = 1.69
$ws.Range("S187").Value = "23/12/2024 23:02"
$ws.Range("T187").Value = 1.65
$ws.Range("U187").Value = "07/01/2024 12:21"
$ws.Range("V187").Value = "https://www.betexplorer.com/football/italy/serie-a/empoli-ac-milan/vPAOOCx2/"

$ws.Range("A188").Value = 187
$ws.Range("B188").Value = "italy"
$ws.Range("C188").Value = "serie-a"
$ws.Range("D188").Value = "2023-2024"
$ws.Range("E188").Value = 45298.625
$ws.Range("F188").Value = "Torino"
$ws.Range("G188").Value = 3
$ws.Range("H188").Value = "Napoli"
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = 3.25
$ws.Range("K188").Value = "23/12/2024 23:02"
$ws.Range("L188").Value = 2.65
$ws.Range("M188").Value = "07/01/2024 14:59"
$ws.Range("N188").Value = 3.31
$ws.Range("O188").Value = "23/12/2024 23:02"
$ws.Range("P188").Value = 3.11
$ws.Range("Q188").Value = "07/01/2024 14:53"
$ws.Range("R188").Value = 2.2
$ws.Range("S188").Value = "23/12/2024 23:02"
$ws.Range("T188").Value = 3.04
$ws.Range("U188").Value = "07/01/2024 14:59"
$ws.Range("V188").Value = "https://www.betexplorer.com/football/italy/serie-a/torino-napoli/0nigvhxk/"

$ws.Range("A189").Value = 188
$ws.Range("B189").Value = "italy"
$ws.Range("C189").Value = "serie-a"
$ws.Range("D189").Value = "2023-2024"
$ws.Range("E189").Value = 45298.625
$ws.Range("F189").Value = "Udinese"
$ws.Range("G189").Value = 1
$ws.Range("H189").Value = "Lazio"
$ws.Range("I189").Value = 2
$ws.Range("J189").Value = 2.95
$ws.Range("K189").Value = "23/12/2024 23:02"
$ws.Range("L189").Value = 3.63
$ws.Range("M189").Value = "07/01/2024 14:59"
$ws.Range("N189").Value = 3.23
$ws.Range("O189").Value = "23/12/2024 23:02"
$ws.Range("P189").Value = 3.23
$ws.Range("Q189").Value = "07/01/2024 14:40"
$ws.Range("R189").Value = 2.4
$ws.Range("S189").Value = "23/12/2024 23:02"
$ws.Range("T189").Value = 2.25
$ws.Range("U189").Value = "07/01/2024 14:53"
$ws.Range("V189").Value = "https://www.betexplorer.com/football/italy/serie-a/udinese-lazio/vwjcwCie/"
